# Refresh the scraped crypto price/volume figures (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.368.16'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '1.843.74'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.08'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6299'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07472'
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2900'
$ws.Range('E9').Value = '  +0.20%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.00'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07728'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.839.88'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.979'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6763'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001031'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.90'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.234'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').Value = '29.369.69'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.60'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.32'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.401'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '158.30'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.509'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1355'
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.49'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06631'
$ws.Range('E28').Value = '  +9.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.439'
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.065'
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.054'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.835'
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.140'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6998'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.582'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01847'
$ws.Range('E37').Value = '  +1.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.824'
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Value = '1.235.87'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.767'
$ws.Range('E40').Value = '  +4.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9360'
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9994'
$ws.Range('D43').Value = '1.993.85'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.11'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.52'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.046'
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.714'
$ws.Range('E48').Value = '  +3.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.957'
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1144'
$ws.Range('E50').Value = '  -1.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3918'
$ws.Range('E51').Value = '  -0.42%  '
